$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F6").Value = 98
$ws.Range("G6").Value = 2928.24
$ws.Range("B10").Value = 32266.79
$ws.Range("F68").Value = 59
$ws.Range("G68").Value = 6792.08
$ws.Range("F78").Value = 45
$ws.Range("G78").Value = 2560.5
$ws.Range("F84").Value = 48
$ws.Range("G84").Value = 4918.08
$ws.Range("B90").Value = 202565.33
$ws.Range("F115").Value = 233
$ws.Range("G115").Value = 22556.73
$ws.Range("B117").Value = 16609.01
$ws.Range("F144").Value = 1234
$ws.Range("G144").Value = 10427.3
$ws.Range("B147").Value = 18779.48
$ws.Range("F149").Value = 254
$ws.Range("G149").Value = 16459.2
$ws.Range("B156").Value = 36230.29
$ws.Range("B219").Value = 63565
$ws.Range("E219").Value = 109.19
$ws.Range("F219").Value = 60
$ws.Range("G219").Value = 6162.6
$ws.Range("B220").Value = 61610
$ws.Range("E220").Value = 122.71
$ws.Range("F220").Value = -58
$ws.Range("G220").Value = -5957.18
$ws.Range("B227").Value = 63520
$ws.Range("E227").Value = 153.4
$ws.Range("F227").Value = 67
$ws.Range("G227").Value = 9666.76
$ws.Range("B228").Value = 55373
$ws.Range("E228").Value = 163.62
$ws.Range("F228").Value = -94
$ws.Range("G228").Value = -13562.32
$ws.Range("B229").Value = 57802
$ws.Range("E229").Value = 162.71
$ws.Range("F229").Value = -79
$ws.Range("G229").Value = -11334.92
$ws.Range("B230").Value = 63531
$ws.Range("E230").Value = 152.53
$ws.Range("F230").Value = 69
$ws.Range("G230").Value = 9900.120000000001
$ws.Range("B232").Value = 55356
$ws.Range("E232").Value = 54.04
$ws.Range("F232").Value = -158
$ws.Range("G232").Value = -7527.12
$ws.Range("B233").Value = 63510
$ws.Range("E233").Value = 50.66
$ws.Range("F233").Value = 127
$ws.Range("G233").Value = 6050.28
$ws.Range("B243").Value = 60325
$ws.Range("E243").Value = 151.57
$ws.Range("F243").Value = -102
$ws.Range("G243").Value = -12939.72
$ws.Range("B244").Value = 63560
$ws.Range("E244").Value = 134.87
$ws.Range("F244").Value = 1
$ws.Range("G244").Value = 126.86
$ws.Range("F255").Value = 611
$ws.Range("G255").Value = 104682.63
$ws.Range("B260").Value = 211085.17
$ws.Range("F270").Value = 44
$ws.Range("G270").Value = 1418.56
$ws.Range("F273").Value = 19
$ws.Range("G273").Value = 6099.19
$ws.Range("F274").Value = 4
$ws.Range("G274").Value = 1284.04
$ws.Range("B275").Value = 7732.76
$ws.Range("F278").Value = 18
$ws.Range("G278").Value = 2471.76
$ws.Range("F282").Value = 12
$ws.Range("G282").Value = 644.4
$ws.Range("F284").Value = 1
$ws.Range("G284").Value = 53.95
$ws.Range("B304").Value = 199483.68
$ws.Range("F328").Value = 64
$ws.Range("G328").Value = 2381.44
$ws.Range("B330").Value = 32445.16
$ws.Range("F354").Value = 22
$ws.Range("G354").Value = 1508.98
$ws.Range("B358").Value = 37597.57
$ws.Range("B375").Value = 45718
$ws.Range("E375").Value = 19.38
$ws.Range("F375").Value = -294
$ws.Range("G375").Value = -4768.68
$ws.Range("B376").Value = 64927
$ws.Range("E376").Value = 17.26
$ws.Range("F376").Value = 106
$ws.Range("G376").Value = 1719.32
$ws.Range("B385").Value = 65067
$ws.Range("E385").Value = 15.65
$ws.Range("F385").Value = 126
$ws.Range("G385").Value = 1855.98
$ws.Range("B386").Value = 53595
$ws.Range("E386").Value = 17.61
$ws.Range("F386").Value = -335
$ws.Range("G386").Value = -4934.55
$ws.Range("F409").Value = 11
$ws.Range("G409").Value = 6423.45
$ws.Range("B411").Value = 9979.110000000001
$ws.Range("F430").Value = 8
$ws.Range("G430").Value = 103.12
$ws.Range("F434").Value = 43
$ws.Range("G434").Value = 1403.52
$ws.Range("B435").Value = 1691.9
$ws.Range("B473").Value = 64830
$ws.Range("E473").Value = 34.9
$ws.Range("F473").Value = 109
$ws.Range("G473").Value = 3578.47
$ws.Range("B474").Value = 60022
$ws.Range("E474").Value = 37.22
$ws.Range("F474").Value = -113
$ws.Range("G474").Value = -3709.79
$ws.Range("F482").Value = 44
$ws.Range("G482").Value = 2607.88
$ws.Range("B488").Value = 33085.9
$ws.Range("F509").Value = 249
$ws.Range("G509").Value = 20014.62
$ws.Range("B510").Value = 26562.84
$ws.Range("F551").Value = 9
$ws.Range("G551").Value = 1288.17
$ws.Range("F554").Value = 8
$ws.Range("G554").Value = 298.24
$ws.Range("B560").Value = 9636.17
$ws.Range("F599").Value = 2080
$ws.Range("G599").Value = 339268.8
$ws.Range("B606").Value = 523988.98
$ws.Range("F613").Value = 150
$ws.Range("G613").Value = 23874
$ws.Range("F617").Value = 3
$ws.Range("G617").Value = 118.59
$ws.Range("B618").Value = 47029.05
$ws.Range("B619").Value = 2003889.27
$ws.Range("B620").Value = 2003889.27
